$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ferramentas AWS")

# Update status values: Cost Explorer (row 2) and CUR (row 4) become "Inativo"
$ws.Range("B2").Value = "Inativo"
$ws.Range("B4").Value = "Inativo"

# Update the active selection to B3 as reflected in the saved workbook
$ws.Activate()
$ws.Range("B3").Select()
